$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Duplicate the two existing data rows (2:3) down into rows 4:5 so the new
# rows inherit the same cell style (s="1") used by the rest of the table,
# then overwrite them with the new team data below.
$ws.Range("A2:P3").Copy() | Out-Null
$ws.Range("A4:P5").Insert(-4121) | Out-Null

# Row 4: Barcelona
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Barcelona"
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 59
$ws.Range("E4").Value = 12
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 12
$ws.Range("N4").Value = 12
$ws.Range("O4").Value = 660
$ws.Range("P4").Value = 613

# Row 5: Paris
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Paris"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 41
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 4
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 12
$ws.Range("N5").Value = 11
$ws.Range("O5").Value = 430
$ws.Range("P5").Value = 387

# Match the row height the new rows ended up with (same as rows 2:3)
$ws.Rows.Item(4).RowHeight = 14.25
$ws.Rows.Item(5).RowHeight = 14.25

# Move the selection like the author's session ended up
$ws.Range("F9").Select() | Out-Null
